$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updates to existing rows (values that changed for "yesterday") ---
$ws.Range("B247").Value = 25
$ws.Range("B268").Value = 46
$ws.Range("B269").Value = 32
$ws.Range("B270").Value = 41
$ws.Range("D270").Value = 3

# --- Add new row 271 (new day of data) ---
# Copy formatting (date number format / style) from the row above first
$ws.Range("A270").Copy() | Out-Null
$ws.Range("A271").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

$ws.Range("A271").Value = 44169
$ws.Range("B271").Value = 24
$ws.Range("C271").Formula = "=B271+C270"
$ws.Range("D271").Value = 2
$ws.Range("E271").Formula = "=D271+E270"
$ws.Range("F271").Formula = "=AVERAGE(B265:B271)"

# --- Update the view: keep header row frozen, update active selection ---
$win = $excel.ActiveWindow
$ws.Range("B2").Select() | Out-Null
$win.FreezePanes = $false
$win.FreezePanes = $true

$ws.Range("L271").Select() | Out-Null
